$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value parses as a plain number need to be forced to
# Text format first, otherwise Excel auto-converts the literal into a
# numeric cell instead of keeping it as the original text string.

$ws.Range('D2').Value = '70.146.21'
$ws.Range('E2').Value = '  -0.40%  '
$ws.Range('D3').Value = '2.518.67'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.96'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.42'
$ws.Range('E6').Value = '  -2.15%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.521'
$ws.Range('E8').Value = '  +1.89%  '
$ws.Range('D9').Value = '2.516.64'
$ws.Range('E9').Value = '  -1.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.162'
$ws.Range('E10').Value = '  -2.23%  '
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.93'
$ws.Range('E13').Value = '  +2.07%  '
$ws.Range('D14').Value = '2.978.59'
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').Value = '69.998.81'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('E16').Value = '  -2.39%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.05'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = '2.516.40'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.45'
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.79'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.67'
$ws.Range('E21').Value = '  -2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.94'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.99'
$ws.Range('E23').Value = '  -1.27%  '
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.57'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.02'
$ws.Range('E26').Value = '  -1.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.88'
$ws.Range('E27').Value = '  -4.95%  '
$ws.Range('D28').Value = '2.656.75'
$ws.Range('E28').Value = '  -1.24%  '
$ws.Range('E29').Value = '  +0.52%  '
$ws.Range('D30').Value = '0.0₃0901'
$ws.Range('E30').Value = '  -3.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.90'
$ws.Range('E31').Value = '  +0.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '465.62'
$ws.Range('E32').Value = '  -4.13%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.26'
$ws.Range('E33').Value = '  -2.45%  '
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '155.92'
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.08'
$ws.Range('E38').Value = '  +1.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.68'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.79'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.320'
$ws.Range('E42').Value = '  -0.66%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.62'
$ws.Range('E43').Value = '  -3.90%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.16'
$ws.Range('E44').Value = '  -13.00%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '38.41'
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.30'
$ws.Range('E46').Value = '  -7.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.88'
$ws.Range('E47').Value = '  -1.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.530'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.50'
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.59'
$ws.Range('E50').Value = '  -3.17%  '
$ws.Range('E51').Value = '  -0.98%  '
